# Project "Sample Project" is saved.
# The Rules worksheet's decision table had its B11 "Rule name" cell
# changed from the text "R40" to the text "1" (still a text value, not a
# number), so the sheet's existing shared-string usage grows by one new
# unique entry "1" while B11 keeps its original cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")

# Writing "1" straight into Value/Formula would be auto-coerced to a
# number by Excel's smart entry, and forcing text via NumberFormat/quote
# -prefix on the cell itself would reformat it (changing its style).
# Instead, stage the literal text "1" in a scratch cell using a formula
# (TEXT() always returns a string), then copy only the *value* (not the
# formatting) into B11 so its existing style/format is left untouched.
$scratch = $ws.Range("Z1048576")
$scratch.Formula = "=TEXT(1,""0"")"

$scratch.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
